$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4: update date and Changes text (Editor/Description columns unchanged)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 45884
$ws.Range("B4").Value = "Changes`n- Made a couple of comments (night of 8/14/25)                                                                                                                                                                                                                                        "

# ---------------------------------------------------------------------------
# Row 5: update date, Changes text, Notes text, and row height
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 45885
$ws.Range("B5").Value = "Changes`n- ALUOp: Changed ALUOp to be the same as IR31to26 aka the OP Code (updated top_level, Datapath, Controller, and alu_control entities)`n- CONTROLLER: Added and defined REG_FETCH, INST_DECODE, and R_TYPE states         `n- ALU CONTROL: Added case for ADD 4 to PC, R-type inst, and ADDU inst                                                                                                                                                                                                                                "
$ws.Range("D5").Value = "Notes`n- Don't think I need to use any kind of states in ALU Control`nBugs`n- Haven't tested any of this yet but it compiles"
$ws.Rows.Item(5).RowHeight = 86.4

# ---------------------------------------------------------------------------
# Rows 6, 7, 8: clear the Date cell only (keep Editor/Description/Notes cells)
# ---------------------------------------------------------------------------
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A8").ClearContents()

# ---------------------------------------------------------------------------
# Rows 9 and 10: brand-new rows appended at the bottom (no Date cell at all)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = $ws.Range("B6").Value()
$ws.Range("B9").WrapText = $true
$ws.Range("B9").VerticalAlignment = -4108

$ws.Range("C9").Value = $ws.Range("C6").Value()
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108

$ws.Range("D9").Value = $ws.Range("D6").Value()
$ws.Range("D9").WrapText = $true
$ws.Range("D9").VerticalAlignment = -4108

$ws.Rows.Item(9).RowHeight = 72

$ws.Range("B10").Value = $ws.Range("B6").Value()
$ws.Range("B10").WrapText = $true
$ws.Range("B10").VerticalAlignment = -4108

$ws.Range("C10").Value = $ws.Range("C6").Value()
$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("C10").VerticalAlignment = -4108

$ws.Range("D10").Value = $ws.Range("D6").Value()
$ws.Range("D10").WrapText = $true
$ws.Range("D10").VerticalAlignment = -4108

$ws.Rows.Item(10).RowHeight = 72

# ---------------------------------------------------------------------------
# Selection: reflects the last-active cell after the edits
# ---------------------------------------------------------------------------
$ws.Range("G9").Select()
